# Generate Report for Handoff
# Refresh the localization-status report: two stale source entries
# (the *.png screenshots that used to carry IsDependency rows) are gone,
# the previous in-flight markdown file (ed1e4cf4-...) has been handed off
# and is replaced by a freshly generated one (da6b6262-...), and a brand
# new markdown file (f317f839-...) has entered the pipeline. The
# ".localization-config" / "Not to be localized" bookkeeping row stays
# last. Net effect: each sheet shrinks from 5 data+header rows to 4.

$wb = $excel.ActiveWorkbook

$e2eBase        = "https://github.com/OpenLocalizationTest/oltest/blob/505f5358432ffa6b594728a7055a2e57b146e36a/e2e/"
$configUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/505f5358432ffa6b594728a7055a2e57b146e36a/.localization-config"
$htZhCnBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e5e9afd7427252bb065bacfcc4d74b2d647fa74/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$htDeDeBase     = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ddcfd60247edbb43aa493097e9129ae961adf1e1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$md1  = "da6b6262-989d-45ef-aa75-f529863dbb53.md"
$md2  = "f317f839-b659-42d9-a200-bd4d2095f5e6.md"

$xlf1ZhCn = "da6b6262-989d-45ef-aa75-f529863dbb53.568f7530e0fb2af58920fc58521c47cd2adb027f.zh-cn.xlf"
$xlf2ZhCn = "f317f839-b659-42d9-a200-bd4d2095f5e6.94aa194dd9071f05b04b2c17b1598eec2744aae1.zh-cn.xlf"
$xlf1DeDe = "da6b6262-989d-45ef-aa75-f529863dbb53.568f7530e0fb2af58920fc58521c47cd2adb027f.de-de.xlf"
$xlf2DeDe = "f317f839-b659-42d9-a200-bd4d2095f5e6.94aa194dd9071f05b04b2c17b1598eec2744aae1.de-de.xlf"

$dtZhCn = "2016-03-09 05:20:54"
$dtDeDe = "2016-03-09 05:20:56"
$epoch  = "0001-01-01 00:00:00"

$ready  = "Ready for handoff"
$notLoc = "Not to be localized"
$config = ".localization-config"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()

$ws1.Range("A2").Value = $md1
$ws1.Range("B2").Value = $ready
$ws1.Range("C2").Value = $ready

$ws1.Range("A3").Value = $md2
$ws1.Range("B3").Value = $ready
$ws1.Range("C3").Value = $ready

$ws1.Range("A4").Value = $config
$ws1.Range("B4").Value = $notLoc
$ws1.Range("C4").Value = $notLoc

$ws1.Rows.Item(5).Delete()

$ws1.Hyperlinks.Add($ws1.Range("A2"), ($e2eBase + $md1), "", "", $md1) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($e2eBase + $md2), "", "", $md2) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, "", "", $config) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

$ws2.Range("A2").Value = $md1
$ws2.Range("B2").Value = $ready
$ws2.Range("C2").Value = $xlf1ZhCn
$ws2.Range("D2").Value = $dtZhCn
$ws2.Range("G2").Value = $epoch
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = $md2
$ws2.Range("B3").Value = $ready
$ws2.Range("C3").Value = $xlf2ZhCn
$ws2.Range("D3").Value = $dtZhCn
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = $config
$ws2.Range("B4").Value = $notLoc
$ws2.Range("C4").Clear()
$ws2.Range("D4").Value = $epoch
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = "Ignored"
$ws2.Range("I2").Clear()
$ws2.Range("I3").Clear()

$ws2.Rows.Item(5).Delete()

$ws2.Hyperlinks.Add($ws2.Range("A2"), ($e2eBase + $md1), "", "", $md1) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), ($htZhCnBase + $xlf1ZhCn), "", "", $xlf1ZhCn) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), ($e2eBase + $md2), "", "", $md2) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), ($htZhCnBase + $xlf2ZhCn), "", "", $xlf2ZhCn) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, "", "", $config) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

$ws3.Range("A2").Value = $md1
$ws3.Range("B2").Value = $ready
$ws3.Range("C2").Value = $xlf1DeDe
$ws3.Range("D2").Value = $dtDeDe
$ws3.Range("G2").Value = $epoch
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = $md2
$ws3.Range("B3").Value = $ready
$ws3.Range("C3").Value = $xlf2DeDe
$ws3.Range("D3").Value = $dtDeDe
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = $config
$ws3.Range("B4").Value = $notLoc
$ws3.Range("C4").Clear()
$ws3.Range("D4").Value = $epoch
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = "Ignored"
$ws3.Range("I2").Clear()
$ws3.Range("I3").Clear()

$ws3.Rows.Item(5).Delete()

$ws3.Hyperlinks.Add($ws3.Range("A2"), ($e2eBase + $md1), "", "", $md1) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), ($htDeDeBase + $xlf1DeDe), "", "", $xlf1DeDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), ($e2eBase + $md2), "", "", $md2) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), ($htDeDeBase + $xlf2DeDe), "", "", $xlf2DeDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, "", "", $config) | Out-Null

Write-Host "Report regenerated for handoff."
